$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50: add work package entry (Controller-AdminControlShowcases, 3h, 2019-01-03)
$ws.Range("D50").Value = "Controller-AdminControlShowcases"
$ws.Range("E50").Value = 3
$ws.Range("G49").Copy($ws.Range("G50"))
$ws.Range("G50").Value = 43468

# Row 51: add work package entry (Model-Datenbank-Operationen, 1h, 2019-01-03)
$ws.Range("D51").Value = "Model-Datenbank-Operationen"
$ws.Range("E51").Value = 1
$ws.Range("G51").Value = 43468

# Update the view selection to match the author's last-saved cursor position
$ws.Range("G52").Select()
